$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.153.11'
$ws.Range("E2").Value = '  +0.95%  '
$ws.Range("D3").Value = '2.420.32'
$ws.Range("E3").Value = '  +1.52%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '563.23'
$ws.Range("E5").Value = '  +1.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.31'
$ws.Range("E6").Value = '  +2.79%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +1.30%  '
$ws.Range("D9").Value = '2.417.07'
$ws.Range("E9").Value = '  +1.30%  '
$ws.Range("E10").Value = '  +0.65%  '
$ws.Range("E11").Value = '  -1.70%  '
$ws.Range("E12").Value = '  +0.81%  '
$ws.Range("E13").Value = '  +0.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.04'
$ws.Range("E14").Value = '  +2.03%  '
$ws.Range("E15").Value = '  +4.73%  '
$ws.Range("D16").Value = '2.856.54'
$ws.Range("E16").Value = '  +2.03%  '
$ws.Range("D17").Value = '62.049.58'
$ws.Range("E17").Value = '  +0.87%  '
$ws.Range("D18").Value = '2.419.08'
$ws.Range("E18").Value = '  +1.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.27'
$ws.Range("E19").Value = '  +2.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '325.11'
$ws.Range("E20").Value = '  +1.14%  '
$ws.Range("E21").Value = '  +0.92%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.77'
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.46'
$ws.Range("E24").Value = '  +1.51%  '
$ws.Range("E25").Value = '  -1.95%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.03'
$ws.Range("E26").Value = '  +0.86%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '587.32'
$ws.Range("E27").Value = '  +12.36%  '
$ws.Range("E28").Value = '  +1.66%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  +3.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.24'
$ws.Range("E31").Value = '  -0.20%  '
$ws.Range("E32").Value = '  +4.08%  '
$ws.Range("E33").Value = '  -0.44%  '
$ws.Range("E34").Value = '  +2.11%  '
$ws.Range("E35").Value = '  +0.22%  '
$ws.Range("E36").Value = '  +2.30%  '
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("E38").Value = '  +0.80%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '154.13'
$ws.Range("E39").Value = '  +4.91%  '
$ws.Range("E40").Value = '  +0.95%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.71'
$ws.Range("E41").Value = '  +0.66%  '
$ws.Range("E42").Value = '  -4.06%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.35'
$ws.Range("E44").Value = '  +7.61%  '
$ws.Range("E45").Value = '  +0.83%  '
$ws.Range("E46").Value = '  +0.89%  '
$ws.Range("E47").Value = '  +2.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.38'
$ws.Range("E48").Value = '  +2.85%  '
$ws.Range("E49").Value = '  +1.33%  '
$ws.Range("E50").Value = '  +1.84%  '
$ws.Range("E51").Value = '  +1.68%  '